$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# Row 2 - TC_Search_1
# ---------------------------------------------------------------------

$ws.Range("C2").Value = @"
*stable internet connection
*chrome browser
"@

$ws.Range("D2").Value = @"

*Open chrome
*Write "Search Key" 
in the search bar 
*Click enter
"@

$ws.Range("E2").Value = 'Search Key="Test Automation"' + [char]10 + 'Search result="What Is Test Automation? A Simple, Clear Introduction"'
$ws.Range("E2").Characters(30, 13).Font.Bold = $true

$ws.Range("F2").Value = "The First Search result " + [char]10 + 'is displayed  ' + [char]10 + '"What Is Test Automation? A Simple, Clear Introduction"'
$ws.Range("F2").Characters(10, 14).Font.Bold = $true

# ---------------------------------------------------------------------
# Row 3 - TC_Search_2
# ---------------------------------------------------------------------

$ws.Range("C3").Value = @"
*stable internet connection
*chrome browser
"@

$ws.Range("D3").Value = @"
*Open chrome
*Write "Search Key" 
in the search bar 
*Click enter
*Scroll down
*Click Next
"@

$ws.Range("E3").Value = 'Search Key="Test Automation"' + [char]10

# ---------------------------------------------------------------------
# Row 4 - TC_Search_3
# ---------------------------------------------------------------------

$ws.Range("C4").Value = @"
*stable internet connection
*chrome browser
"@

$ws.Range("D4").Value = @"
*Open chrome
*Write "Search Key" 
in the search bar 
*Click enter
*Scroll down
*Click Next
*Scroll down
*Click Previous
"@

$ws.Range("E4").Value = 'Search Key="Test Automation"' + [char]10

# ---------------------------------------------------------------------
# Final selection, matches the recorded cursor position after editing
# ---------------------------------------------------------------------
$ws.Range("D4").Select() | Out-Null
